# Fix problems with graph and graph generator spreadsheet
#
# The June 2023 circulation count (row 7, column F) was actually July's
# figure. Correct June's value and fill in July's real value, which had
# been left blank in column F (row 8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F7").Value = 114493
$ws.Range("F8").Value = 101818
